$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$values = @{
    "B2" = -0.3291291692495283;  "C2" = 1.812906911408057;  "D2" = 27.04923173616717;  "E2" = 5.200887591187409;  "F2" = 5.247817881751421;  "G2" = 46
    "B3" = -0.9495326209681835;  "C3" = 2.023872778234117;  "D3" = 43.75773521907136;  "E3" = 6.614962979418053;  "F3" = 6.620432449373355;  "G3" = 45
    "B4" = -1.78880049929551;    "C4" = 3.145904294961889;  "D4" = 171.4597706064012;  "E4" = 13.094264798239;    "F4" = 13.12147041243131;  "G4" = 44
    "B5" = -4.213288228997171;   "C5" = 5.762860054837126;  "D5" = 837.7929586707833;  "E5" = 28.94465336933202;  "F5" = 28.97526445982997;  "G5" = 43
    "B6" = -8.901962788184429;   "C6" = 10.44294745540641;  "D6" = 3576.380529151133;  "E6" = 59.80284716592625;  "F6" = 59.85341742634296;  "G6" = 42
    "B7" = -18.90835081736334;   "C7" = 20.50494998719962;  "D7" = 14831.91797591911;  "E7" = 121.786362027606;   "F7" = 121.8041544916215;  "G7" = 41
    "B8" = -39.82066791606469;   "C8" = 41.45267336896583;  "D8" = 64161.89026247886;  "E8" = 253.3019744543632;  "F8" = 253.3391456874329;  "G8" = 40
    "B9" = -84.10188452280053;   "C9" = 85.72607652365316;  "D9" = 276745.0078080686;  "E9" = 526.0655926859963;  "F9" = 526.087929219649;   "G9" = 39
    "B10" = -176.9409422996708;  "C10" = 178.6784928052071; "D10" = 1193366.26590965;  "E10" = 1092.413047299258; "F10" = 1092.458304185648; "G10" = 38
    "B11" = -372.8822845639866;  "C11" = 374.5397916143799; "D11" = 5147954.627584639; "E11" = 2268.910449441458; "F11" = 2268.931447530792; "G11" = 37
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
